# Update figures in Sheet1 for the 2022-05-10 Fonds de solidarite data refresh.
# Only column C (nombre_aides), column D (nombre_entreprises) and column E
# (montant_total) change for a subset of rows, per the published diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 766296
$ws.Range("E2").Value = 1429156227

$ws.Range("C13").Value = 187830
$ws.Range("E13").Value = 1165045094

$ws.Range("C54").Value = 75189
$ws.Range("E54").Value = 361025801

$ws.Range("C69").Value = 17884
$ws.Range("E69").Value = 103656076

$ws.Range("C78").Value = 178439
$ws.Range("E78").Value = 892467990

$ws.Range("C88").Value = 71246
$ws.Range("E88").Value = 110277178

$ws.Range("C91").Value = 18836
$ws.Range("E91").Value = 75024804

$ws.Range("C93").Value = 16911
$ws.Range("E93").Value = 50425027

$ws.Range("C100").Value = 9326
$ws.Range("E100").Value = 23686296

$ws.Range("C104").Value = 319348
$ws.Range("E104").Value = 561269974

$ws.Range("C112").Value = 145222
$ws.Range("E112").Value = 715925957

$ws.Range("C115").Value = 81797
$ws.Range("E115").Value = 436221192

$ws.Range("C121").Value = 1306072
$ws.Range("D121").Value = 220381
$ws.Range("E121").Value = 2274394878

$ws.Range("C129").Value = 633247
$ws.Range("E129").Value = 3425036750

$ws.Range("C130").Value = 4238
$ws.Range("E130").Value = 140270462

$ws.Range("C132").Value = 585527
$ws.Range("E132").Value = 3457024600

$ws.Range("C136").Value = 26669
$ws.Range("D136").Value = 4272
$ws.Range("E136").Value = 143465662

$ws.Range("C139").Value = 76627
$ws.Range("E139").Value = 114124865

$ws.Range("C144").Value = 25047
$ws.Range("E144").Value = 92100271

$ws.Range("C151").Value = 39911
$ws.Range("E151").Value = 60347099

$ws.Range("C154").Value = 18423
$ws.Range("E154").Value = 72343317

$ws.Range("C156").Value = 12389
$ws.Range("E156").Value = 39983570

$ws.Range("C164").Value = 196201
$ws.Range("E164").Value = 370288687

$ws.Range("C177").Value = 6937
$ws.Range("E177").Value = 30788203

$ws.Range("C178").Value = 515872
$ws.Range("E178").Value = 891180160

$ws.Range("C186").Value = 236810
$ws.Range("E186").Value = 1189625057

$ws.Range("C196").Value = 595489
$ws.Range("E196").Value = 983987711

$ws.Range("C221").Value = 135491
$ws.Range("D221").Value = 27175
$ws.Range("E221").Value = 681804100

$ws.Range("C237").Value = 283301
$ws.Range("E237").Value = 1438274452

$ws.Range("C246").Value = 18836
$ws.Range("E246").Value = 71609583
